$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.003.24"
$ws.Range("E2").Value = "  -0.81%  "
$ws.Range("D3").Value = "1.619.95"
$ws.Range("E3").Value = "  -1.39%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.57"
$ws.Range("E5").Value = "  -1.21%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.518"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.251"
$ws.Range("E8").Value = "  -1.28%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0627"
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.09"
$ws.Range("E10").Value = "  +0.55%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0848"
$ws.Range("E11").Value = "  +0.01%  "
$ws.Range("D12").Value = "1.848.14"
$ws.Range("E12").Value = "  -1.38%  "
$ws.Range("D13").Value = "1.606.34"
$ws.Range("E13").Value = "  -2.47%  "
$ws.Range("E14").Value = "  +0.04%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.539"
$ws.Range("E15").Value = "  -0.68%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.72"
$ws.Range("E16").Value = "  -3.55%  "
$ws.Range("D17").Value = "26.986.54"
$ws.Range("E17").Value = "  -0.86%  "
$ws.Range("D18").Value = "0.0₃0747"
$ws.Range("E18").Value = "  +0.77%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "213.68"
$ws.Range("E19").Value = "  -2.48%  "
$ws.Range("E20").Value = "  -0.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.85"
$ws.Range("E21").Value = "  -1.88%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.36"
$ws.Range("E22").Value = "  -1.28%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.40"
$ws.Range("E23").Value = "  -5.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.05"
$ws.Range("E24").Value = "  -0.90%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "148.11"
$ws.Range("E25").Value = "  +0.05%  "
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.38"
$ws.Range("E27").Value = "  -1.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.116"
$ws.Range("E28").Value = "  -2.01%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.55"
$ws.Range("E29").Value = "  -1.09%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0514"
$ws.Range("E31").Value = "  -1.15%  "
$ws.Range("E32").Value = "  -1.12%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.749"
$ws.Range("E33").Value = "  +35.84%  "
$ws.Range("E34").Value = "  -0.12%  "
$ws.Range("D35").Value = "1.342.85"
$ws.Range("E35").Value = "  +2.61%  "
$ws.Range("E36").Value = "  -0.61%  "
$ws.Range("E37").Value = "  -0.51%  "
$ws.Range("E38").Value = "  +0.41%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.847"
$ws.Range("E39").Value = "  -1.30%  "
$ws.Range("E40").Value = "  -0.12%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.801"
$ws.Range("E41").Value = "  -1.35%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.22"
$ws.Range("E42").Value = "  -0.28%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "65.11"
$ws.Range("E43").Value = "  +5.22%  "
$ws.Range("E44").Value = "  -0.08%  "
$ws.Range("D45").Value = "1.758.80"
$ws.Range("E45").Value = "  -1.43%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "89.66"
$ws.Range("E46").Value = "  -2.39%  "
$ws.Range("B47").Value = "WEMIXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.876"
$ws.Range("E47").Value = "  +30.94%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.62"
$ws.Range("E48").Value = "  +1.43%  "
$ws.Range("E49").Value = "  -0.25%  "
$ws.Range("E50").Value = "  +4.42%  "
$ws.Range("E51").Value = "  +0.51%  "
